$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column AB (old AB..AL shift right to AD..AN)
$ws.Range("AB1:AC1").EntireColumn.Insert()

# New column headers
$ws.Range("AB1").Value = "terminal compressed hydrogen storage amount (days)"
$ws.Range("AC1").Value = "terminal liquid hydrogen storage amount (days)"

# New baseline values (row 2 only - the rest of the rows stay blank in these columns)
$ws.Range("AB2").Value = 0.25
$ws.Range("AC2").Value = 1

# Move the two comments that were anchored on the shifted cells (old AJ1/AK1)
# so they stay attached to the same logical column headers (now AL1/AM1).
$commentAJ = $ws.Range("AJ1").Comment
$textAJ = $commentAJ.Text()
$commentAJ.Delete()
$ws.Range("AL1").AddComment($textAJ)

$commentAK = $ws.Range("AK1").Comment
$textAK = $commentAK.Text()
$commentAK.Delete()
$ws.Range("AM1").AddComment($textAK)

# Extend the conditional formatting range to cover the two new columns
$fc = $ws.Range("C3:AL56").FormatConditions
$rule = $fc.Item(1)
$rule.ModifyAppliesToRange($ws.Range("C3:AN56"))

# Update selection / active cell to match the post-edit view
$ws.Range("AB1").Select()
